$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subtotals")

# 1. Sort the list (A2:E43) by Product (column B), ascending, with header row.
$dataRange = $ws.Range("A1:E43")
$dataRange.Sort($ws.Range("B1"), 1, $null, $null, 1, $null, $null, 1)

# 2. Apply Data > Subtotal: at each change in Product (col 2), use Sum (function 9),
#    add subtotal to Sales (col 5), replacing current subtotals, summary below data off,
#    grand total on.
$dataRange2 = $ws.Range("A1:E43")
$dataRange2.Subtotal(2, 9, @(5), $true, $false, $true)

# 3. Collapse the outline: detail rows go to outline level 2 and are hidden,
#    the per-product subtotal rows stay at outline level 1 (visible, collapsed group).
$detailGroups = @(
    @(2, 8),
    @(10, 18),
    @(20, 29),
    @(31, 46)
)
foreach ($grp in $detailGroups) {
    $first = $grp[0]
    $last = $grp[1]
    $rows = $ws.Range("A$($first):A$($last)").EntireRow
    $rows.OutlineLevel = 2
    $rows.Hidden = $true
}

# 4. Activate the Subtotals sheet and select the cell below the new data, matching
#    the end-state selection/tab captured in the workbook.
$ws.Activate()
$ws.Range("B50").Select()
